# Weekly update: insert 3 new "Espárragos" price records for
# Vega Modelo de Temuco at the top of the data block (rows 75-77),
# pushing the existing historical rows down by three positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new rows right before the current row 75.
$ws.Range("A75:A77").EntireRow.Insert()

# --- New row 75 ---
$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "Vega Modelo de Temuco"
$ws.Range("C75").Value = "La Araucanía"
$ws.Range("D75").Value2 = 45218
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = 300000000
$ws.Range("G75").Value = "Espárragos"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Extra"
$ws.Range("J75").Value = 185
$ws.Range("K75").Value = 2000
$ws.Range("L75").Value = 2000
$ws.Range("M75").Value = 2000
$ws.Range("N75").Value = "$/kilo"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 2000
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"

# --- New row 76 ---
$ws.Range("A76").Value = 10
$ws.Range("B76").Value = "Vega Modelo de Temuco"
$ws.Range("C76").Value = "La Araucanía"
$ws.Range("D76").Value2 = 45218
$ws.Range("E76").Value = 9
$ws.Range("F76").Value = 300000000
$ws.Range("G76").Value = "Espárragos"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 500
$ws.Range("K76").Value = 1800
$ws.Range("L76").Value = 1800
$ws.Range("M76").Value = 1800
$ws.Range("N76").Value = "$/kilo"
$ws.Range("O76").Value = "Región de La Araucanía"
$ws.Range("P76").Value = 1800
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# --- New row 77 ---
$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "Vega Modelo de Temuco"
$ws.Range("C77").Value = "La Araucanía"
$ws.Range("D77").Value2 = 45218
$ws.Range("E77").Value = 9
$ws.Range("F77").Value = 300000000
$ws.Range("G77").Value = "Espárragos"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 600
$ws.Range("K77").Value = 1500
$ws.Range("L77").Value = 1500
$ws.Range("M77").Value = 1500
$ws.Range("N77").Value = "$/kilo"
$ws.Range("O77").Value = "Región del Maule"
$ws.Range("P77").Value = 1500
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"
